$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Faiblesses column (B3:B5) - refreshed wording
$ws.Range("B3").Value = "Equipe dynamique"
$ws.Range("B4").Value = "Equipe toujours en veille technologique"
$ws.Range("B5").Value = "Efficacité de réaction face à l'imprévu "

# Menaces column (E3:E4) - refreshed wording / new entry
$ws.Range("E3").Value = "Manque d'experience de projet "
$ws.Range("E4").Value = "Connaissances fonctionnelles du BTP"

# Opportunités column (B9:B10) - new entries
$ws.Range("B9").Value = "Augmenter l'image de marque de la société"
$ws.Range("B10").Value = 'Acquerir de l"experience dans le domaine'

# Menaces (bottom block, E10) - new entry
$ws.Range("E10").Value = "Concurrence et futur évolutions "

# Update the selection to match the saved view state
$ws.Range("E5:G5").Select()
